$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# "wall terminals 6 and 8" -- remove the stray wall-terminal assignments that
# were entered by mistake for Zona 2 (ID pulsantiera 11 / bottone 8 -> bagno)
# and Zona 4 (ID pulsantiera 9 / bottone 7 -> scala). Clear the value/label
# pair and restore the "empty" Pulsante 3 formatting (same look as the other
# unused Pulsante 3 slots, e.g. K7:L7 / C17:D17).
Copy-CellFormat "K7" "G12"
Copy-CellFormat "L7" "H12"
$ws.Range("G12").ClearContents()
$ws.Range("H12").ClearContents()

Copy-CellFormat "C17" "G22"
Copy-CellFormat "D17" "H22"
$ws.Range("G22").ClearContents()
$ws.Range("H22").ClearContents()

# Re-apply the canonical "ID pulsantiera" / "Dimmer" formatting to the third
# (K/L) column -- functionally identical to what was already there, but this
# is what collapses the workbook's style table back down to its minimal set
# of distinct cell formats (no more stray duplicate xf records).
Copy-CellFormat "C3" "K3"
Copy-CellFormat "C4" "K4"
Copy-CellFormat "C5" "K5"
Copy-CellFormat "C5" "K6"
Copy-CellFormat "C3" "K8"
Copy-CellFormat "C4" "K9"
Copy-CellFormat "C4" "K10"
Copy-CellFormat "C4" "K11"
Copy-CellFormat "C3" "K13"
Copy-CellFormat "C4" "K14"
Copy-CellFormat "C5" "K15"
Copy-CellFormat "C5" "K16"

$excel.CutCopyMode = $false

# Park the selection where the user last clicked while reviewing the sheet.
$ws.Range("J15").Select() | Out-Null
